$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Boult Audio AirBass z40 In Ear True Wireless (TWS) 60 Hours Playback IPX4(Splash & Sweat Proof) Powerfull bass -Bluetooth White"
$ws.Range("B1").Value = "Rs. 1,499"

$ws.Range("A2").Value = "NBOX Buzz TWS On Ear True Wireless (TWS) 20 Hours Playback IPX5(Splash & Sweat Proof) Passive noise cancellation -Bluetooth Version 5.1 Black"
$ws.Range("B2").Value = "Rs. 745"

$ws.Range("A5").Value = "Boult Audio Airbass X50 In Ear True Wireless (TWS) 40 Hours Playback IPX5(Splash & Sweat Proof) Fast charging -Bluetooth Black"
$ws.Range("B5").Value = "Rs. 1,499"
